$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2Project")

$ws.Range("A2").Value = "2202263747"
$ws.Range("B2").Value = "Automation Project7666947"
